$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C, rows 2 through 43 all currently hold the serial date 45849
# (2025-07-11) and need to be bumped to 45850 (2025-07-12).
for ($row = 2; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45849) {
        $cell.Value2 = 45850
    }
}
